# Scheduled runner update: refresh market price snapshots (currentAveragePrice,
# currentAveragePriceNQ/HQ) and the dependent Leve price / profit columns
# (LevePriceNQ/HQ, LeveProfitNQ/HQ) across the ALC, ARM, BSM, CRP, CUL, GSM,
# LTW and WVR sheets. A few rows saw their market data drop to zero, in which
# case the corresponding profit cell is cleared entirely rather than left at
# a stale value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 231.66667  # H12: 235.4 -> 231.66667
$ws.Cells.Item(12, 9).Value = 197.5  # I12: 199 -> 197.5
$ws.Cells.Item(12, 10).Value = 300  # J12: 290 -> 300
$ws.Cells.Item(12, 11).Value = 197.5  # K12: 199 -> 197.5
$ws.Cells.Item(12, 12).Value = 300  # L12: 290 -> 300
$ws.Cells.Item(12, 13).Value = -27.5  # M12: -29 -> -27.5
$ws.Cells.Item(12, 14).Value = -640  # N12: -630 -> -640

$ws.Cells.Item(33, 8).Value = 132  # H33: 136.38461 -> 132
$ws.Cells.Item(33, 9).Value = 130.3077  # I33: 134.91667 -> 130.3077
$ws.Cells.Item(33, 11).Value = 130.3077  # K33: 134.91667 -> 130.3077
$ws.Cells.Item(33, 13).Value = 98.69229999999999  # M33: 94.08332999999999 -> 98.69229999999999

$ws.Cells.Item(70, 8).Value = 4746.5713  # H70: 5291.9165 -> 4746.5713
$ws.Cells.Item(70, 10).Value = 5204.3335  # J70: 5950.3 -> 5204.3335
$ws.Cells.Item(70, 12).Value = 15613.0005  # L70: 17850.9 -> 15613.0005
$ws.Cells.Item(70, 14).Value = -16153.0005  # N70: -18390.9 -> -16153.0005

$ws.Cells.Item(73, 8).Value = 4746.5713  # H73: 5291.9165 -> 4746.5713
$ws.Cells.Item(73, 10).Value = 5204.3335  # J73: 5950.3 -> 5204.3335
$ws.Cells.Item(73, 12).Value = 15613.0005  # L73: 17850.9 -> 15613.0005
$ws.Cells.Item(73, 14).Value = -17485.0005  # N73: -19722.9 -> -17485.0005

$ws.Cells.Item(105, 8).Value = 67478.336  # H105: 68717.5 -> 67478.336
$ws.Cells.Item(105, 10).Value = 66217.5  # J105: 67435 -> 66217.5
$ws.Cells.Item(105, 12).Value = 66217.5  # L105: 67435 -> 66217.5
$ws.Cells.Item(105, 14).Value = -73205.5  # N105: -74423 -> -73205.5

$ws.Cells.Item(108, 8).Value = 0  # H108: 61000 -> 0
$ws.Cells.Item(108, 10).Value = 0  # J108: 61000 -> 0
$ws.Cells.Item(108, 12).Value = 0  # L108: 61000 -> 0
$ws.Cells.Item(108, 14).ClearContents()  # N108: -68680 -> (removed)

$ws.Cells.Item(135, 8).Value = 1178.5454  # H135: 1059 -> 1178.5454
$ws.Cells.Item(135, 9).Value = 1061.9  # I135: 951.8333 -> 1061.9
$ws.Cells.Item(135, 11).Value = 9557.1  # K135: 8566.4997 -> 9557.1
$ws.Cells.Item(135, 13).Value = -7022.1  # M135: -6031.4997 -> -7022.1

$ws.Cells.Item(138, 8).Value = 4720.2427  # H138: 4721.6714 -> 4720.2427
$ws.Cells.Item(138, 10).Value = 5963.684  # J138: 5966.316 -> 5963.684
$ws.Cells.Item(138, 12).Value = 17891.052  # L138: 17898.948 -> 17891.052
$ws.Cells.Item(138, 14).Value = -28171.052  # N138: -28178.948 -> -28171.052

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1816.7142  # H2: 1610.5555 -> 1816.7142
$ws.Cells.Item(2, 9).Value = 1426.75  # I2: 1319.2 -> 1426.75
$ws.Cells.Item(2, 10).Value = 2336.6667  # J2: 1974.75 -> 2336.6667
$ws.Cells.Item(2, 11).Value = 1426.75  # K2: 1319.2 -> 1426.75
$ws.Cells.Item(2, 12).Value = 2336.6667  # L2: 1974.75 -> 2336.6667
$ws.Cells.Item(2, 13).Value = -1313.75  # M2: -1206.2 -> -1313.75
$ws.Cells.Item(2, 14).Value = -2562.6667  # N2: -2200.75 -> -2562.6667

$ws.Cells.Item(32, 8).Value = 4144.784  # H32: 4066.423 -> 4144.784
$ws.Cells.Item(32, 9).Value = 2586.422  # I32: 2531.7173 -> 2586.422
$ws.Cells.Item(32, 11).Value = 2586.422  # K32: 2531.7173 -> 2586.422
$ws.Cells.Item(32, 13).Value = -2299.422  # M32: -2244.7173 -> -2299.422

$ws.Cells.Item(61, 8).Value = 1950  # H61: 1983.3334 -> 1950
$ws.Cells.Item(61, 9).Value = 1950  # I61: 1983.3334 -> 1950
$ws.Cells.Item(61, 11).Value = 1950  # K61: 1983.3334 -> 1950
$ws.Cells.Item(61, 13).Value = -1738  # M61: -1771.3334 -> -1738

$ws.Cells.Item(116, 8).Value = 1816.7142  # H116: 1610.5555 -> 1816.7142
$ws.Cells.Item(116, 9).Value = 1426.75  # I116: 1319.2 -> 1426.75
$ws.Cells.Item(116, 10).Value = 2336.6667  # J116: 1974.75 -> 2336.6667
$ws.Cells.Item(116, 11).Value = 1426.75  # K116: 1319.2 -> 1426.75
$ws.Cells.Item(116, 12).Value = 2336.6667  # L116: 1974.75 -> 2336.6667
$ws.Cells.Item(116, 13).Value = 867.25  # M116: 974.8 -> 867.25
$ws.Cells.Item(116, 14).Value = -6924.6667  # N116: -6562.75 -> -6924.6667

$ws.Cells.Item(122, 8).Value = 5318.8335  # H122: 5215.409 -> 5318.8335
$ws.Cells.Item(122, 9).Value = 4436.5835  # I122: 4514.9375 -> 4436.5835
$ws.Cells.Item(122, 11).Value = 13309.7505  # K122: 13544.8125 -> 13309.7505
$ws.Cells.Item(122, 13).Value = -10859.7505  # M122: -11094.8125 -> -10859.7505

$ws.Cells.Item(136, 8).Value = 1950  # H136: 1983.3334 -> 1950
$ws.Cells.Item(136, 9).Value = 1950  # I136: 1983.3334 -> 1950
$ws.Cells.Item(136, 11).Value = 5850  # K136: 5950.0002 -> 5850
$ws.Cells.Item(136, 13).Value = -3300  # M136: -3400.0002 -> -3300

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1816.7142  # H3: 1610.5555 -> 1816.7142
$ws.Cells.Item(3, 9).Value = 1426.75  # I3: 1319.2 -> 1426.75
$ws.Cells.Item(3, 10).Value = 2336.6667  # J3: 1974.75 -> 2336.6667
$ws.Cells.Item(3, 11).Value = 1426.75  # K3: 1319.2 -> 1426.75
$ws.Cells.Item(3, 12).Value = 2336.6667  # L3: 1974.75 -> 2336.6667
$ws.Cells.Item(3, 13).Value = -1312.75  # M3: -1205.2 -> -1312.75
$ws.Cells.Item(3, 14).Value = -2564.6667  # N3: -2202.75 -> -2564.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(69, 8).Value = 0  # H69: 25000 -> 0
$ws.Cells.Item(69, 9).Value = 0  # I69: 25000 -> 0
$ws.Cells.Item(69, 11).Value = 0  # K69: 25000 -> 0
$ws.Cells.Item(69, 13).ClearContents()  # M69: -24251 -> (removed)

$ws.Cells.Item(72, 8).Value = 0  # H72: 25000 -> 0
$ws.Cells.Item(72, 9).Value = 0  # I72: 25000 -> 0
$ws.Cells.Item(72, 11).Value = 0  # K72: 75000 -> 0
$ws.Cells.Item(72, 13).ClearContents()  # M72: -71256 -> (removed)

$ws.Cells.Item(86, 8).Value = 6200  # H86: 6333.3335 -> 6200
$ws.Cells.Item(86, 9).Value = 4000  # I86: 4200 -> 4000
$ws.Cells.Item(86, 10).Value = 15000  # J86: 17000 -> 15000
$ws.Cells.Item(86, 11).Value = 4000  # K86: 4200 -> 4000
$ws.Cells.Item(86, 12).Value = 15000  # L86: 17000 -> 15000
$ws.Cells.Item(86, 13).Value = -2877  # M86: -3077 -> -2877
$ws.Cells.Item(86, 14).Value = -17246  # N86: -19246 -> -17246

$ws.Cells.Item(89, 8).Value = 6200  # H89: 6333.3335 -> 6200
$ws.Cells.Item(89, 9).Value = 4000  # I89: 4200 -> 4000
$ws.Cells.Item(89, 10).Value = 15000  # J89: 17000 -> 15000
$ws.Cells.Item(89, 11).Value = 20000  # K89: 21000 -> 20000
$ws.Cells.Item(89, 12).Value = 75000  # L89: 85000 -> 75000
$ws.Cells.Item(89, 13).Value = -14384  # M89: -15384 -> -14384
$ws.Cells.Item(89, 14).Value = -86232  # N89: -96232 -> -86232

$ws.Cells.Item(134, 8).Value = 3786.4546  # H134: 3971.5 -> 3786.4546
$ws.Cells.Item(134, 9).Value = 3739.2222  # I134: 3962.3333 -> 3739.2222
$ws.Cells.Item(134, 11).Value = 11217.6666  # K134: 11886.9999 -> 11217.6666
$ws.Cells.Item(134, 13).Value = -8682.6666  # M134: -9351.999899999999 -> -8682.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(62, 8).Value = 4874.875  # H62: 3656.25 -> 4874.875
$ws.Cells.Item(62, 10).Value = 3666.5  # J62: 2964.2856 -> 3666.5
$ws.Cells.Item(62, 12).Value = 10999.5  # L62: 8892.856800000001 -> 10999.5
$ws.Cells.Item(62, 14).Value = -12371.5  # N62: -10264.8568 -> -12371.5

$ws.Cells.Item(65, 8).Value = 4874.875  # H65: 3656.25 -> 4874.875
$ws.Cells.Item(65, 10).Value = 3666.5  # J65: 2964.2856 -> 3666.5
$ws.Cells.Item(65, 12).Value = 32998.5  # L65: 26678.5704 -> 32998.5
$ws.Cells.Item(65, 14).Value = -39862.5  # N65: -33542.5704 -> -39862.5

$ws.Cells.Item(109, 8).Value = 892.3333  # H109: 1001.125 -> 892.3333
$ws.Cells.Item(109, 9).Value = 892.3333  # I109: 1001.125 -> 892.3333
$ws.Cells.Item(109, 11).Value = 2676.9999  # K109: 3003.375 -> 2676.9999
$ws.Cells.Item(109, 13).Value = -1636.9999  # M109: -1963.375 -> -1636.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5567.3335  # H70: 5786.5713 -> 5567.3335
$ws.Cells.Item(70, 9).Value = 4390.1113  # I70: 4626.625 -> 4390.1113
$ws.Cells.Item(70, 11).Value = 4390.1113  # K70: 4626.625 -> 4390.1113
$ws.Cells.Item(70, 13).Value = -4120.1113  # M70: -4356.625 -> -4120.1113

$ws.Cells.Item(73, 8).Value = 5567.3335  # H73: 5786.5713 -> 5567.3335
$ws.Cells.Item(73, 9).Value = 4390.1113  # I73: 4626.625 -> 4390.1113
$ws.Cells.Item(73, 11).Value = 4390.1113  # K73: 4626.625 -> 4390.1113
$ws.Cells.Item(73, 13).Value = -3454.1113  # M73: -3690.625 -> -3454.1113

$ws.Cells.Item(80, 8).Value = 15389.8  # H80: 12138.23 -> 15389.8
$ws.Cells.Item(80, 9).Value = 7499.75  # I80: 5499.6665 -> 7499.75
$ws.Cells.Item(80, 10).Value = 20649.834  # J80: 17828.428 -> 20649.834
$ws.Cells.Item(80, 11).Value = 7499.75  # K80: 5499.6665 -> 7499.75
$ws.Cells.Item(80, 12).Value = 20649.834  # L80: 17828.428 -> 20649.834
$ws.Cells.Item(80, 13).Value = -6501.75  # M80: -4501.6665 -> -6501.75
$ws.Cells.Item(80, 14).Value = -22645.834  # N80: -19824.428 -> -22645.834

$ws.Cells.Item(83, 8).Value = 15389.8  # H83: 12138.23 -> 15389.8
$ws.Cells.Item(83, 9).Value = 7499.75  # I83: 5499.6665 -> 7499.75
$ws.Cells.Item(83, 10).Value = 20649.834  # J83: 17828.428 -> 20649.834
$ws.Cells.Item(83, 11).Value = 37498.75  # K83: 27498.3325 -> 37498.75
$ws.Cells.Item(83, 12).Value = 103249.17  # L83: 89142.14 -> 103249.17
$ws.Cells.Item(83, 13).Value = -32506.75  # M83: -22506.3325 -> -32506.75
$ws.Cells.Item(83, 14).Value = -113233.17  # N83: -99126.14 -> -113233.17

$ws.Cells.Item(93, 8).Value = 59326.668  # H93: 25768.908 -> 59326.668
$ws.Cells.Item(93, 10).Value = 59326.668  # J93: 25768.908 -> 59326.668
$ws.Cells.Item(93, 12).Value = 59326.668  # L93: 25768.908 -> 59326.668
$ws.Cells.Item(93, 14).Value = -63070.668  # N93: -29512.908 -> -63070.668

$ws.Cells.Item(122, 8).Value = 95511.91  # H122: 87886.086 -> 95511.91
$ws.Cells.Item(122, 9).Value = 4237.5713  # I122: 4444.1665 -> 4237.5713
$ws.Cells.Item(122, 10).Value = 255242  # J122: 171328 -> 255242
$ws.Cells.Item(122, 11).Value = 12712.7139  # K122: 13332.4995 -> 12712.7139
$ws.Cells.Item(122, 12).Value = 765726  # L122: 513984 -> 765726
$ws.Cells.Item(122, 13).Value = -10262.7139  # M122: -10882.4995 -> -10262.7139
$ws.Cells.Item(122, 14).Value = -770626  # N122: -518884 -> -770626

$ws.Cells.Item(123, 8).Value = 46199.9  # H123: 31500.334 -> 46199.9
$ws.Cells.Item(123, 10).Value = 46199.9  # J123: 31500.334 -> 46199.9
$ws.Cells.Item(123, 12).Value = 46199.9  # L123: 31500.334 -> 46199.9
$ws.Cells.Item(123, 14).Value = -51099.9  # N123: -36400.334 -> -51099.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 5419.7144  # H16: 5419.857 -> 5419.7144
$ws.Cells.Item(16, 9).Value = 5823  # I16: 5823.1665 -> 5823
$ws.Cells.Item(16, 11).Value = 5823  # K16: 5823.1665 -> 5823
$ws.Cells.Item(16, 13).Value = -5653  # M16: -5653.1665 -> -5653

$ws.Cells.Item(22, 8).Value = 2319.8865  # H22: 2279.4443 -> 2319.8865
$ws.Cells.Item(22, 10).Value = 3671.6428  # J22: 3460.2 -> 3671.6428
$ws.Cells.Item(22, 12).Value = 3671.6428  # L22: 3460.2 -> 3671.6428
$ws.Cells.Item(22, 14).Value = -4261.6428  # N22: -4050.2 -> -4261.6428

$ws.Cells.Item(27, 8).Value = 2319.8865  # H27: 2279.4443 -> 2319.8865
$ws.Cells.Item(27, 10).Value = 3671.6428  # J27: 3460.2 -> 3671.6428
$ws.Cells.Item(27, 12).Value = 3671.6428  # L27: 3460.2 -> 3671.6428
$ws.Cells.Item(27, 14).Value = -3885.6428  # N27: -3674.2 -> -3885.6428

$ws.Cells.Item(36, 8).Value = 0  # H36: 79888 -> 0
$ws.Cells.Item(36, 10).Value = 0  # J36: 79888 -> 0
$ws.Cells.Item(36, 12).Value = 0  # L36: 79888 -> 0
$ws.Cells.Item(36, 14).ClearContents()  # N36: -81012 -> (removed)

$ws.Cells.Item(68, 8).Value = 4400.857  # H68: 4350.75 -> 4400.857
$ws.Cells.Item(68, 9).Value = 4266.6665  # I68: 4200 -> 4266.6665
$ws.Cells.Item(68, 11).Value = 4266.6665  # K68: 4200 -> 4266.6665
$ws.Cells.Item(68, 13).Value = -3517.6665  # M68: -3451 -> -3517.6665

$ws.Cells.Item(71, 8).Value = 4400.857  # H71: 4350.75 -> 4400.857
$ws.Cells.Item(71, 9).Value = 4266.6665  # I71: 4200 -> 4266.6665
$ws.Cells.Item(71, 11).Value = 21333.3325  # K71: 21000 -> 21333.3325
$ws.Cells.Item(71, 13).Value = -17589.3325  # M71: -17256 -> -17589.3325

$ws.Cells.Item(76, 8).Value = 11500  # H76: 34500 -> 11500
$ws.Cells.Item(76, 10).Value = 11500  # J76: 34500 -> 11500
$ws.Cells.Item(76, 12).Value = 11500  # L76: 34500 -> 11500
$ws.Cells.Item(76, 14).Value = -12176  # N76: -35176 -> -12176

$ws.Cells.Item(79, 8).Value = 11500  # H79: 34500 -> 11500
$ws.Cells.Item(79, 10).Value = 11500  # J79: 34500 -> 11500
$ws.Cells.Item(79, 12).Value = 11500  # L79: 34500 -> 11500
$ws.Cells.Item(79, 14).Value = -13840  # N79: -36840 -> -13840

$ws.Cells.Item(82, 8).Value = 2577  # H82: 3034 -> 2577
$ws.Cells.Item(82, 9).Value = 2577  # I82: 3795.6667 -> 2577
$ws.Cells.Item(82, 10).Value = 0  # J82: 749 -> 0
$ws.Cells.Item(82, 11).Value = 2577  # K82: 3795.6667 -> 2577
$ws.Cells.Item(82, 12).Value = 0  # L82: 749 -> 0
$ws.Cells.Item(82, 13).Value = -2216  # M82: -3434.6667 -> -2216
$ws.Cells.Item(82, 14).ClearContents()  # N82: -1471 -> (removed)

$ws.Cells.Item(85, 8).Value = 2577  # H85: 3034 -> 2577
$ws.Cells.Item(85, 9).Value = 2577  # I85: 3795.6667 -> 2577
$ws.Cells.Item(85, 10).Value = 0  # J85: 749 -> 0
$ws.Cells.Item(85, 11).Value = 2577  # K85: 3795.6667 -> 2577
$ws.Cells.Item(85, 12).Value = 0  # L85: 749 -> 0
$ws.Cells.Item(85, 13).Value = -1329  # M85: -2547.6667 -> -1329
$ws.Cells.Item(85, 14).ClearContents()  # N85: -3245 -> (removed)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(80, 8).Value = 33333.332  # H80: 65000 -> 33333.332
$ws.Cells.Item(80, 10).Value = 33333.332  # J80: 65000 -> 33333.332
$ws.Cells.Item(80, 12).Value = 33333.332  # L80: 65000 -> 33333.332
$ws.Cells.Item(80, 14).Value = -35329.332  # N80: -66996 -> -35329.332

$ws.Cells.Item(83, 8).Value = 33333.332  # H83: 65000 -> 33333.332
$ws.Cells.Item(83, 10).Value = 33333.332  # J83: 65000 -> 33333.332
$ws.Cells.Item(83, 12).Value = 99999.99600000001  # L83: 195000 -> 99999.99600000001
$ws.Cells.Item(83, 14).Value = -109983.996  # N83: -204984 -> -109983.996

$ws.Cells.Item(126, 8).Value = 3252.353  # H126: 3093.7368 -> 3252.353
$ws.Cells.Item(126, 9).Value = 3168.3333  # I126: 2901.5 -> 3168.3333
$ws.Cells.Item(126, 10).Value = 3346.875  # J126: 3307.3333 -> 3346.875
$ws.Cells.Item(126, 11).Value = 9504.999899999999  # K126: 8704.5 -> 9504.999899999999
$ws.Cells.Item(126, 12).Value = 10040.625  # L126: 9921.999899999999 -> 10040.625
$ws.Cells.Item(126, 13).Value = -7034.999899999999  # M126: -6234.5 -> -7034.999899999999
$ws.Cells.Item(126, 14).Value = -14980.625  # N126: -14861.9999 -> -14980.625

$ws.Cells.Item(132, 8).Value = 107311.664  # H132: 137832.14 -> 107311.664
$ws.Cells.Item(132, 9).Value = 118227  # I132: 157472.67 -> 118227
$ws.Cells.Item(132, 11).Value = 354681  # K132: 472418.01 -> 354681
$ws.Cells.Item(132, 13).Value = -352151  # M132: -469888.01 -> -352151
